# UC011 - Listar Solicitações de Diárias workbook update
# From v1.1.1 to v1.2
#
# The second test step of TC2 ("Chefe Clica para ordenar pelo nome do
# servidor." / "SYSTEM Visualiza os registros de solicitações de diária
# ordenado pelo nome do servidor.") and the second test step of TC4
# ("Chefe Indica alguns parâmetros específicos para a busca; Informa o
# nome do beneficiário; Filtra a listagem de solicitações." / "SYSTEM
# Exibe uma nova listagem de solicitações, de acordo com os filtros
# informados pelo usuário.") are swapped with each other, so that TC2
# now holds the filter scenario and TC4 now holds the sort scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# TC2 block, second step (row 20): columns B (Steps) and D (Expected Results)
$tc2StepOld = $ws.Range("B20").Value2
$tc2ExpectedOld = $ws.Range("D20").Value2

# TC4 block, second step (row 36): columns B (Steps) and D (Expected Results)
$tc4StepOld = $ws.Range("B36").Value2
$tc4ExpectedOld = $ws.Range("D36").Value2

# Swap the contents
$ws.Range("B20").Value = $tc4StepOld
$ws.Range("D20").Value = $tc4ExpectedOld

$ws.Range("B36").Value = $tc2StepOld
$ws.Range("D36").Value = $tc2ExpectedOld
